# TopSky Developer Guide Settings workbook update
# - Incorporate TopSky 2.3: drop the two now-unused blank sheets and rename
#   the remaining data sheet to "TopSky 2.3"
# - Reset the saved view back to the top of the sheet (frozen header still
#   at row 2, scrolled/selected near the top instead of deep in the list)
# - Protect the worksheet (objects/scenarios locked along with the sheet)

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Drop the empty Sheet2 / Sheet3 tabs -----------------------------------
$wb.Worksheets("Sheet2").Delete() | Out-Null
$wb.Worksheets("Sheet3").Delete() | Out-Null

# --- Rename the remaining sheet --------------------------------------------
$ws = $wb.Worksheets("Sheet1")
$ws.Name = "TopSky 2.3"

# --- Restore the view to the top of the frozen data ------------------------
$ws.Activate() | Out-Null
$ws.Range("A11").Select() | Out-Null

# --- Protect the sheet (objects + scenarios locked, as well as the sheet) --
$ws.Protect($null, $true, $true, $true) | Out-Null
